# Update the "Förändrad" (Changed) date column for all data rows
# from 2023-09-02 (serial 45171) to 2023-09-03 (serial 45172).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45172
}
